$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values look numeric (contain dots as thousands separators
# or decimals) but must remain plain text, exactly as authored upstream.
# Force text format before assignment so Excel does not reinterpret them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.851.64"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.219.65"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.26"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.41"
$ws.Range("E6").Value = "  +5.64%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.39"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.47"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.41"
$ws.Range("E14").Value = "  +1.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.563.23"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.10"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.218.32"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.730"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.810.26"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("E20").Value = "  +11.42%  "

$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.83"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.77"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.75"
$ws.Range("E24").Value = "  +2.80%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.82"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -2.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.26"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.75"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.68"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  +2.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0719"
$ws.Range("E35").Value = "  +2.65%  "

$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("E37").Value = "  +5.82%  "

$ws.Range("E38").Value = "  +1.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.93"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0992"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("E41").Value = "  +2.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.096.13"
$ws.Range("E42").Value = "  +9.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.80"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("E44").Value = "  +6.44%  "

$ws.Range("E45").Value = "  +2.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.00"
$ws.Range("E46").Value = "  +8.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.75"
$ws.Range("E47").Value = "  +7.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.435.97"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.87"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.30"
$ws.Range("E51").Value = "  +0.95%  "
